$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename variable names to shorter GUI-friendly identifiers
$ws.Range("B3").Value = "T_rad"
$ws.Range("B6").Value = "K_f"
$ws.Range("B7").Value = "K_r"

# Update the active selection to reflect where the user clicked next (below the table)
$ws.Range("B9").Select()
